$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row labels
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# Fix capitalization of "Ciudad de México" -> "Ciudad De México"
$ws.Range("A9").Value = "Ciudad De México"

# Remove the footer rows 25-29 (sample size, source, author, dept, date notes)
$ws.Range("A25:D29").EntireRow.Delete()
